# Bugfixing and Pyromancy content
# - adds a new "Burn" worksheet (burn-damage-over-time calculator) after
#   the existing "Abilities and values" sheet and makes it the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# New worksheet, placed right after the existing one.
$burn = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$burn.Name = "Burn"

# ---- Header row (row 1) --------------------------------------------------
# Filled in this exact order so the shared-string table matches the
# original authoring order (Burn I..Burn X, then Resistance, Steps, Damage).
$burn.Range("D1").Value = "Burn I"
$burn.Range("E1").Value = "Burn II"
$burn.Range("F1").Value = "Burn III"
$burn.Range("G1").Value = "Burn IV"
$burn.Range("H1").Value = "Burn V"
$burn.Range("I1").Value = "Burn VI"
$burn.Range("J1").Value = "Burn VII"
$burn.Range("K1").Value = "Burn VIII"
$burn.Range("L1").Value = "Burn IX"
$burn.Range("M1").Value = "Burn X"
$burn.Range("B1").Value = "Resistance"
$burn.Range("C1").Value = "Steps"
$burn.Range("A1").Value = "Damage"

# ---- Step-multiplier row (row 2) ----------------------------------------
$burn.Range("D2").Value = 0.1
$burn.Range("E2").Value = 0.2
$burn.Range("F2").Value = 0.3
$burn.Range("G2").Value = 0.4
$burn.Range("H2").Value = 0.5
$burn.Range("I2").Value = 0.6
$burn.Range("J2").Value = 0.7
$burn.Range("K2").Value = 0.8
$burn.Range("L2").Value = 0.9
$burn.Range("M2").Value = 1

# ---- Data / formula rows -------------------------------------------------
$cols = @("D","E","F","G","H","I","J","K","L","M")

function Set-BurnFormulaRow($sheet, [int]$row) {
    foreach ($col in $cols) {
        $f = "=IF(`$A$row*" + $col + "`$2-`$B$row<=0,0,(`$A$row*" + $col + "`$2-`$B$row)+IF(0.5*`$A$row*" + $col + "`$2-1.5*`$B$row<=0,0,0.5*`$A$row*" + $col + "`$2-1.5*`$B$row)+IF(0.25*`$A$row*" + $col + "`$2-1.75*`$B$row<=0,0,0.25*`$A$row*" + $col + "`$2-1.75*`$B$row+IF(0.125*`$A$row*" + $col + "`$2-1.875*`$B$row<=0,0,0.125*`$A$row*" + $col + "`$2-1.875*`$B$row+IF(0.0625*`$A$row*" + $col + "`$2-1.9375*`$B$row<=0,0,0.0625*`$A$row*" + $col + "`$2-1.9375*`$B$row))))"
        $sheet.Range($col + $row).Formula = $f
    }
}

# 100-value block (rows 4-10), B stepping from 0 to 30 by 5
$aValue = 100
$row = 4
foreach ($bValue in 0,5,10,15,20,25,30) {
    $burn.Range("A$row").Value = $aValue
    $burn.Range("B$row").Value = $bValue
    Set-BurnFormulaRow $burn $row
    $row++
}

# 50-value block (rows 13-19), B stepping from 0 to 30 by 5
$aValue = 50
$row = 13
foreach ($bValue in 0,5,10,15,20,25,30) {
    $burn.Range("A$row").Value = $aValue
    $burn.Range("B$row").Value = $bValue
    Set-BurnFormulaRow $burn $row
    $row++
}

# Column B only needs to fit "Resistance" / the numeric step values.
$burn.Columns.Item(2).AutoFit()

# Selection on the new sheet, matching the author's last selection.
$burn.Range("G4:M4").Select()

# The new sheet becomes the active / selected tab.
$burn.Activate()
